# Fruta / hortaliza, semanal
# Update weekly price/volume data for rows 2-5 (Feria Lagunitas de Puerto Montt - Tuna)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44253
$ws.Range("M2").Value = 160

# Row 3
$ws.Range("D3").Value = 44250
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 806

# Row 4
$ws.Range("D4").Value = 44257
$ws.Range("M4").Value = 100

# Row 5
$ws.Range("D5").Value = 44252
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 13500
$ws.Range("S5").Value = 750
